# "Updated symbol list ... with GitHub Actions" -- refresh crypto prices,
# rotate three token rows, and tweak a couple of Worst/Best-in-24h label
# suffixes, as captured in the commit diff for cryptos.xlsx.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($cellRef, $value) {
    # Some of the new prices keep a trailing zero (e.g. 0.8125 -> 0.8140) or
    # would otherwise be re-rendered in scientific notation by plain numeric
    # auto-detection (e.g. 0.00002099). Formatting the cell as Text first
    # forces Excel to keep the literal digits exactly as provided, matching
    # the original inline-string content.
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
}

# --- Column D (Price) updates ---
$ws.Range("D2").Value  = "247.56"
$ws.Range("D3").Value  = "21.75"
$ws.Range("D4").Value  = "5.341"
$ws.Range("D5").Value  = "0.05631"
$ws.Range("D6").Value  = "3.428"
$ws.Range("D7").Value  = "6.369"
Set-TextValue "D8" "0.8140"
$ws.Range("D9").Value  = "0.9346"
$ws.Range("D10").Value = "0.1434"
$ws.Range("D11").Value = "0.07497"
$ws.Range("D12").Value = "0.03218"
$ws.Range("D14").Value = "0.09296"
$ws.Range("D15").Value = "3.584"
$ws.Range("D16").Value = "0.001606"
$ws.Range("D17").Value = "0.04717"
Set-TextValue "D18" "0.0005780"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("D19").Value = "0.006278"
$ws.Range("D20").Value = "0.005052"
$ws.Range("D23").Value = "3.763"
Set-TextValue "D28" "0.0003000"
$ws.Range("D40").Value = "0.03946"

# --- Rows 41-43: tokens rotate down one slot
#     (KickToken -> row41, BKEXToken -> row42, CEJI -> row43) ---
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "0.006819"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "0.1063"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "0.003399"
$ws.Range("E43").Value = "42CEJICEJI"

# --- Remaining price / label corrections ---
$ws.Range("D44").Value = "0.008798"
Set-TextValue "D45" "0.00005581"
Set-TextValue "D47" "0.0005500"
$ws.Range("E47").Value = "46ACDXExchangeACXT"
Set-TextValue "D48" "0.7800"
Set-TextValue "D50" "0.00002099"
Set-TextValue "D51" "0.01010"
